$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 25; this shifts the former rows 25-47
# (Z08_B06_P01 ... Z17_B03_P01) down to rows 26-48.
$ws.Rows.Item(25).Insert()

# The freshly inserted row doesn't carry the surrounding data-row style,
# so copy the formatting from the row right below it (the row that used
# to be row 25) back onto the new row 25.
$ws.Range("A26:D26").Copy()
$ws.Range("A25:D25").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# Populate the new row 25 with the new bullet point.
$ws.Range("A25").Value = "Z08_B05_P02"
$ws.Range("B25").Value = "Z08_B05"
$ws.Range("C25").Value = "Soziale Gerechtigkeit steigern, gerechte und gute Löhne stärken"
$ws.Range("D25").Value = "XXXSoziale Gerechtigkeit steigern, gerechte und gute Löhne stärken"
